# Table 2 (sheet3.xml): insert a new "Province-level class" column (D)
# that flags which socioeconomic variables were converted to categorical
# for the province-level analysis. Existing column D (Details) shifts to E.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table 2")
$ws.Activate()

# Insert a new blank column before the current column D ("Details"),
# pushing the old D -> E automatically (values, formatting and all).
$ws.Columns.Item(4).Insert()

# Size the new column similarly to how the author left it.
$ws.Columns.Item(4).ColumnWidth = 21.71

# Header for the new column.
$ws.Range("D1").Value = "Province-level class"

# Rows that mark the predictor as having been converted to a categorical
# variable for the province-level class analysis.
$categoricalRows = @(3, 5, 6, 7, 8, 9, 10, 13, 14, 15, 16, 17, 18, 19)
foreach ($r in $categoricalRows) {
    $cell = $ws.Range("D$r")
    $cell.Style = "Normal"
    $cell.Value = "Categorical"
}

# Row 2 (Total population) was not converted -> marked "NA".
$ws.Range("D2").Style = "Normal"
$ws.Range("D2").Value = "NA"

# Restore the selection the author left the sheet on.
$ws.Range("D4").Select()
